$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "UNLOCK"
$ws.Range("F6").Value = "UNLOCK"
$ws.Range("E10").Value = "NEUTRAL"
$ws.Range("F10").Value = "NEUTRAL"

$ws.Range("I13").Select()
